$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.924.71'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.783.69'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.544'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.18%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.99'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0682'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0936'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = '2.040.61'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.96%  '
$ws.Range('D14').Value = '1.743.16'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('D15').Value = '33.918.08'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.615'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '0.0₃0772'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0513'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.62'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').Value = '1.396.05'
$ws.Range('E35').Value = '  +0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.643'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.33'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.02%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.78'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.35%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.36'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.918'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.67'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +15.32%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0139'
$ws.Range('E45').Value = '  +14.09%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.53%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0508'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.71'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('D50').Value = '1.942.56'
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('E51').Value = '  -0.06%  '
